$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.338.88"
$ws.Range("E2").Value = "  -1.90%  "

# Row 3
$ws.Range("D3").Value = "1.853.70"
$ws.Range("E3").Value = "  -1.17%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7007"
$ws.Range("E5").Value = "  -5.65%  "

# Row 6
$ws.Range("E6").Value = "  -1.28%  "

# Row 7
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3082"
$ws.Range("E8").Value = "  -2.34%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07490"
$ws.Range("E9").Value = "  +3.72%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.83"
$ws.Range("E10").Value = "  -3.46%  "

# Row 11
$ws.Range("E11").Value = "  -3.36%  "

# Row 12
$ws.Range("D12").Value = "1.879.70"
$ws.Range("E12").Value = "  +0.29%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7282"
$ws.Range("E13").Value = "  -3.10%  "

# Row 14
$ws.Range("E14").Value = "  -3.87%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.49"
$ws.Range("E15").Value = "  -3.36%  "

# Row 16
$ws.Range("D16").Value = "29.546.24"
$ws.Range("E16").Value = "  -1.20%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.917"
$ws.Range("E17").Value = "  -2.73%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "242.92"
$ws.Range("E18").Value = "  -1.91%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007758"
$ws.Range("E19").Value = "  -1.20%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.17"
$ws.Range("E20").Value = "  -2.99%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  +0.29%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").Value = "2.143.00"
$ws.Range("E22").Value = "  +0.58%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.0000"
$ws.Range("E23").Value = "  -0.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.648"
$ws.Range("E24").Value = "  -4.83%  "

# Row 25
$ws.Range("E25").Value = "  -5.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.054"
$ws.Range("E26").Value = "  -2.27%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.97"
$ws.Range("E27").Value = "  -1.88%  "

# Row 28
$ws.Range("E28").Value = "  -2.78%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.947"
$ws.Range("E29").Value = "  -4.53%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.389"

# Row 31
$ws.Range("E31").Value = "  -1.61%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.414"
$ws.Range("E32").Value = "  -4.38%  "

# Row 33
$ws.Range("E33").Value = "  -4.85%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05278"
$ws.Range("E34").Value = "  -0.79%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.203"
$ws.Range("E35").Value = "  -2.81%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7235"
$ws.Range("E36").Value = "  -3.97%  "

# Row 37
$ws.Range("E37").Value = "  +0.45%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.679"
$ws.Range("E38").Value = "  -0.38%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01872"
$ws.Range("E39").Value = "  -4.57%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.716"
$ws.Range("E40").Value = "  -1.54%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8884"
$ws.Range("E41").Value = "  +4.05%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4321"
$ws.Range("E42").Value = "  -4.48%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.928"
$ws.Range("E43").Value = "  -1.92%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "70.28"
$ws.Range("E44").Value = "  -3.18%  "

# Row 45
$ws.Range("D45").Value = "1.051.37"
$ws.Range("E45").Value = "  -5.51%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  -0.06%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.67"
$ws.Range("E47").Value = "  -0.85%  "

# Row 48
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.270"
$ws.Range("E48").Value = "  -4.66%  "

# Row 49
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.034.92"
$ws.Range("E49").Value = "  +0.57%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.758"
$ws.Range("E50").Value = "  -5.39%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.299"
$ws.Range("E51").Value = "  -1.59%  "
